$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 45, shifting existing rows 45-143 down to 46-144
$ws.Rows(45).Insert()

# Populate the newly inserted row 45 with the new weekly data point
$ws.Range("A45").Value = 4
$ws.Range("B45").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C45").Value = "Los Lagos"
$ws.Range("D45").Value = 44526
$ws.Range("E45").Value = 10
$ws.Range("F45").Value = 100112039
$ws.Range("G45").Value = "Ciboulette"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 240
$ws.Range("K45").Value = 2500
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = 2500
$ws.Range("N45").Value = "`$/docena de atados"
$ws.Range("O45").Value = "Región Metropolitana"
$ws.Range("P45").Value = 833
$ws.Range("Q45").Value = 3
$ws.Range("R45").Value = "Hortaliza"
